$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column C. This shifts the existing
# "Additional Feedback" column (C) to column D, and leaves a blank
# column C ready for the new "Rating" column.
$ws.Columns("C:C").Insert()

# Remove the now-obsolete rows 3 through 10 (only the header row and
# the first data row remain afterwards).
$ws.Rows("3:10").Delete()

# Fill in the new "Rating" header and its data, plus the feedback text.
$ws.Range("C1").Value = "Rating"
$ws.Range("C2").Value = "Maybe"
$ws.Range("D2").Value = "Idk about this person"
